# Update vm_pu.xlsx values for "case with 380 kV done"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bf = New-Object 'object[,]' 24,5
$bf[0,0] = 1.02
$bf[0,1] = 1.075020570840446
$bf[0,2] = 1.074961202190801
$bf[0,3] = 1.078745159909038
$bf[0,4] = 1.084366223436723
$bf[1,0] = 1.02
$bf[1,1] = 1.076591953828826
$bf[1,2] = 1.07625163041827
$bf[1,3] = 1.080269297306914
$bf[1,4] = 1.085864528931957
$bf[2,0] = 1.02
$bf[2,1] = 1.077605632912954
$bf[2,2] = 1.077083359009521
$bf[2,3] = 1.081252735489641
$bf[2,4] = 1.086831119334081
$bf[3,0] = 1.02
$bf[3,1] = 1.078031051717314
$bf[3,2] = 1.077432246681844
$bf[3,3] = 1.081665519086292
$bf[3,4] = 1.087236788209337
$bf[4,0] = 1.02
$bf[4,1] = 1.078102438859624
$bf[4,2] = 1.077490781548832
$bf[4,3] = 1.08173478927141
$bf[4,4] = 1.087304861926508
$bf[5,0] = 1.02
$bf[5,1] = 1.077611320237993
$bf[5,2] = 1.077088023878745
$bf[5,3] = 1.081258253676977
$bf[5,4] = 1.086836542579724
$bf[6,0] = 1.02
$bf[6,1] = 1.075552278953304
$bf[6,2] = 1.075397990829849
$bf[6,3] = 1.079260832205263
$bf[6,4] = 1.084873193148427
$bf[7,0] = 1.02
$bf[7,1] = 1.071899580885731
$bf[7,2] = 1.072394455793066
$bf[7,3] = 1.07571928843943
$bf[7,4] = 1.081390676818579
$bf[8,0] = 1.02
$bf[8,1] = 1.069447226151252
$bf[8,2] = 1.070374296232096
$bf[8,3] = 1.073342827201148
$bf[8,4] = 1.079052894324078
$bf[9,0] = 1.02
$bf[9,1] = 1.068381061915519
$bf[9,2] = 1.069495171671327
$bf[9,3] = 1.072309964002858
$bf[9,4] = 1.078036621496489
$bf[10,0] = 1.02
$bf[10,1] = 1.067984382411376
$bf[10,2] = 1.069167954205956
$bf[10,3] = 1.071925721061676
$bf[10,4] = 1.077658517392308
$bf[11,0] = 1.02
$bf[11,1] = 1.06806950157264
$bf[11,2] = 1.069238174059939
$bf[11,3] = 1.072008169483742
$bf[11,4] = 1.07773965007713
$bf[12,0] = 1.02
$bf[12,1] = 1.068348285792264
$bf[12,2] = 1.069468137547452
$bf[12,3] = 1.072278214514645
$bf[12,4] = 1.078005379926373
$bf[13,0] = 1.02
$bf[13,1] = 1.06851996619813
$bf[13,2] = 1.069609736272465
$bf[13,3] = 1.07244451934126
$bf[13,4] = 1.078169022908704
$bf[14,0] = 1.02
$bf[14,1] = 1.069517892469744
$bf[14,2] = 1.070432547393169
$bf[14,3] = 1.073411292781167
$bf[14,4] = 1.079120255543443
$bf[15,0] = 1.02
$bf[15,1] = 1.070142708661836
$bf[15,2] = 1.070947492652593
$bf[15,3] = 1.074016685888731
$bf[15,4] = 1.079715858170128
$bf[16,0] = 1.02
$bf[16,1] = 1.070506741646845
$bf[16,2] = 1.071247429589844
$bf[16,3] = 1.074369431961094
$bf[16,4] = 1.080062878387621
$bf[17,0] = 1.02
$bf[17,1] = 1.070630798127576
$bf[17,2] = 1.071349629163093
$bf[17,3] = 1.074489647026276
$bf[17,4] = 1.080181138496281
$bf[18,0] = 1.02
$bf[18,1] = 1.070075714499055
$bf[18,2] = 1.070892287635372
$bf[18,3] = 1.073951771261158
$bf[18,4] = 1.079651995543115
$bf[19,0] = 1.02
$bf[19,1] = 1.068266209097859
$bf[19,2] = 1.069400437673594
$bf[19,3] = 1.07219870943243
$bf[19,4] = 1.077927146206997
$bf[20,0] = 1.02
$bf[20,1] = 1.067124682585299
$bf[20,2] = 1.068458560929424
$bf[20,3] = 1.071093059898072
$bf[20,4] = 1.076839098714646
$bf[21,0] = 1.02
$bf[21,1] = 1.06773019468974
$bf[21,2] = 1.068958240889557
$bf[21,3] = 1.071679515702283
$bf[21,4] = 1.077416236213441
$bf[22,0] = 1.02
$bf[22,1] = 1.070105987549536
$bf[22,2] = 1.070917233709978
$bf[22,3] = 1.073981104528706
$bf[22,4] = 1.079680853503965
$bf[23,0] = 1.02
$bf[23,1] = 1.072846866233748
$bf[23,2] = 1.073174029810988
$bf[23,3] = 1.076637526641898
$bf[23,4] = 1.082293774162079

$inArr = New-Object 'object[,]' 24,6
$inArr[0,0] = 1.061115250838073
$inArr[0,1] = 1.079927127846741
$inArr[0,2] = 1.077649031706224
$inArr[0,3] = 1.081423027142735
$inArr[0,4] = 1.087029430281817
$inArr[0,5] = 1.081460748237657
$inArr[1,0] = 1.061736971784661
$inArr[1,1] = 1.081154903621809
$inArr[1,2] = 1.078755524041825
$inArr[1,3] = 1.082763385880615
$inArr[1,4] = 1.088345092336323
$inArr[1,5] = 1.082690267595153
$inArr[2,0] = 1.062136163780126
$inArr[2,1] = 1.081945828916188
$inArr[2,2] = 1.079467683090691
$inArr[2,3] = 1.083627401408637
$inArr[2,4] = 1.089192987356396
$inArr[2,5] = 1.08348231609417
$inArr[3,0] = 1.062303247335076
$inArr[3,1] = 1.08227749992742
$inArr[3,2] = 1.079766171322576
$inArr[3,3] = 1.083989857905898
$inArr[3,4] = 1.089548633170818
$inArr[3,5] = 1.083814458116281
$inArr[4,0] = 1.062331258371648
$inArr[4,1] = 1.082333140349558
$inArr[4,2] = 1.07981623612125
$inArr[4,3] = 1.084050670783691
$inArr[4,4] = 1.089608300531097
$inArr[4,5] = 1.083870177554201
$inArr[5,0] = 1.06213839924531
$inArr[5,1] = 1.081950263982927
$inArr[5,2] = 1.079471675041173
$inArr[5,3] = 1.083632247602535
$inArr[5,4] = 1.089197742677835
$inArr[5,5] = 1.083486757459213
$inArr[6,0] = 1.061326010542392
$inArr[6,1] = 1.08034279724616
$inArr[6,2] = 1.078023772040517
$inArr[6,3] = 1.081876694791309
$inArr[6,4] = 1.087474780531151
$inArr[6,5] = 1.0818770079353
$inArr[7,0] = 1.059870436376268
$inArr[7,1] = 1.077482740581487
$inArr[7,2] = 1.075442708551884
$inArr[7,3] = 1.07875751752813
$inArr[7,4] = 1.084411962779209
$inArr[7,5] = 1.07901288966211
$inArr[8,0] = 1.058883500789839
$inArr[8,1] = 1.075556879336268
$inArr[8,2] = 1.073701413602798
$inArr[8,3] = 1.076660103442437
$inArr[8,4] = 1.082351407741242
$inArr[8,5] = 1.077084293473025
$inArr[9,0] = 1.058452133376091
$inArr[9,1] = 1.074718265935972
$inArr[9,2] = 1.072942388425837
$inArr[9,3] = 1.075747481539917
$inArr[9,4] = 1.081454577991577
$inArr[9,5] = 1.076244489145509
$inArr[10,0] = 1.058291293071807
$inArr[10,1] = 1.074406047956521
$inArr[10,2] = 1.072659684132349
$inArr[10,3] = 1.075407814180199
$inArr[10,4] = 1.081120751386261
$inArr[10,5] = 1.075931827780716
$inArr[11,0] = 1.058325821646163
$inArr[11,1] = 1.074473052549173
$inArr[11,2] = 1.072720360141998
$inArr[11,3] = 1.075480704927997
$inArr[11,4] = 1.081192390409171
$inArr[11,5] = 1.075998927527575
$inArr[12,0] = 1.058438850779841
$inArr[12,1] = 1.074692472656737
$inArr[12,2] = 1.07291903575113
$inArr[12,3] = 1.075719418460389
$inArr[12,4] = 1.081426998223219
$inArr[12,5] = 1.076218659236859
$inArr[13,0] = 1.058508410534979
$inArr[13,1] = 1.074827568875616
$inArr[13,2] = 1.07304134413762
$inArr[13,3] = 1.075866407365193
$inArr[13,4] = 1.081571454129103
$inArr[13,5] = 1.076353947307865
$inArr[14,0] = 1.058912043926485
$inArr[14,1] = 1.075612435150389
$inArr[14,2] = 1.073751680509993
$inArr[14,3] = 1.076720576696407
$inArr[14,4] = 1.082410829389289
$inArr[14,5] = 1.077139928182773
$inArr[15,0] = 1.059164151466624
$inArr[15,1] = 1.076103492773755
$inArr[15,2] = 1.074195899571645
$inArr[15,3] = 1.077255179471328
$inArr[15,4] = 1.082936107306778
$inArr[15,5] = 1.07763168316429
$inArr[16,0] = 1.059310814389785
$inArr[16,1] = 1.076389465519967
$inArr[16,2] = 1.07445452027096
$inArr[16,3] = 1.077566577946075
$inArr[16,4] = 1.083242050528348
$inArr[16,5] = 1.077918062024603
$inArr[17,0] = 1.059360757219751
$inArr[17,1] = 1.076486898417244
$inArr[17,2] = 1.074542621394627
$inArr[17,3] = 1.077672684887259
$inArr[17,4] = 1.083346294614868
$inArr[17,5] = 1.07801563328777
$inArr[18,0] = 1.059137142812107
$inArr[18,1] = 1.076050853873711
$inArr[18,2] = 1.074148289301267
$inArr[18,3] = 1.077197865844161
$inArr[18,4] = 1.082879795825759
$inArr[18,5] = 1.07757896951097
$inArr[19,0] = 1.058405583451961
$inArr[19,1] = 1.074627878874942
$inArr[19,2] = 1.072860552052077
$inArr[19,3] = 1.075649142128305
$inArr[19,4] = 1.081357931616473
$inArr[19,5] = 1.076153973724485
$inArr[20,0] = 1.057942082762557
$inArr[20,1] = 1.073729027903331
$inArr[20,2] = 1.072046447708126
$inArr[20,3] = 1.074671462942512
$inArr[20,4] = 1.080396994581039
$inArr[20,5] = 1.075253846281395
$inArr[21,0] = 1.058188131390794
$inArr[21,1] = 1.074205925497069
$inArr[21,2] = 1.072478446078208
$inArr[21,3] = 1.075190126806229
$inArr[21,4] = 1.08090679687277
$inArr[21,5] = 1.075731421124417
$inArr[22,0] = 1.059149348056661
$inArr[22,1] = 1.076074640525992
$inArr[22,2] = 1.074169803830839
$inArr[22,3] = 1.077223764722286
$inArr[22,4] = 1.0829052419273
$inArr[22,5] = 1.077602789943025
$inArr[23,0] = 1.060249625959734
$inArr[23,1] = 1.078225460657449
$inArr[23,2] = 1.07611355673243
$inArr[23,3] = 1.079567013525837
$inArr[23,4] = 1.085207014705405
$inArr[23,5] = 1.079756664485771

$ws.Range("B2:F25").Value = $bf
$ws.Range("I2:N25").Value = $inArr
